{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n// the \"\u00a9 2020 . Contact: ...\" paragraph, and the now-superfluous blank\n// paragraph that separated them from the preceding \"LOB1038: ...\" line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs that must be removed outright.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The empty paragraph immediately before the \"Ver no Jupiter\" line is the\n  // one collapsed away by the edit (two blank separators around the removed\n  // block become one).\n  let blankIdx = -1;\n  for (let i = jupiterIdx - 1; i >= 0; i--) {\n    if (items[i].text === \"\") {\n      blankIdx = i;\n      break;\n    }\n  }\n\n  // Delete from the bottom up so earlier indices stay valid.\n  const toDelete = [jupiterIdx, copyrightIdx];\n  if (blankIdx !== -1) {\n    toDelete.push(blankIdx);\n  }\n  toDelete.sort((a, b) => b - a);\n  for (const idx of toDelete) {\n    items[idx].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n# the \"\u00a9 2020 . Contact: ...\" paragraph, and the now-superfluous blank\n# paragraph that separated them from the preceding \"LOB1038: ...\" line.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\n$copyrightIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -ne -1 -and $copyrightIdx -ne -1) {\n    # The empty paragraph immediately before the \"Ver no Jupiter\" line is the\n    # one collapsed away by the edit (two blank separators around the\n    # removed block become one).\n    $blankIdx = -1\n    for ($i = $jupiterIdx - 1; $i -ge 1; $i--) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        $trimmed = $t.Trim([char]13)\n        if ($trimmed -eq \"\") {\n            $blankIdx = $i\n            break\n        }\n    }\n\n    # Delete from the bottom up so earlier indices stay valid.\n    $toDelete = @($jupiterIdx, $copyrightIdx)\n    if ($blankIdx -ne -1) {\n        $toDelete += $blankIdx\n    }\n    $toDelete = $toDelete | Sort-Object -Descending\n    foreach ($idx in $toDelete) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
